$d = $word.ActiveDocument
$d.Content.Find.Execute("94.6%", $true, $false, $false, $false, $false, $true, 1, $false, "96.4%", 2)
